# DOCS folder for EX2
# The Algo2 comparison table on the right (columns G:J) had its "Lat" (H)
# and "Alt" (J) columns swapped by mistake. Swap the raw H/J values back
# for the two data rows; the H4/J4 "Deviation" formulas (H2-H3 / J2-J3)
# recalc automatically once the inputs change.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tmp = $ws.Range("H2").Value2
$ws.Range("H2").Value2 = $ws.Range("J2").Value2
$ws.Range("J2").Value2 = $tmp

$tmp = $ws.Range("H3").Value2
$ws.Range("H3").Value2 = $ws.Range("J3").Value2
$ws.Range("J3").Value2 = $tmp

# Move the sheet's active cell/selection to J15.
$ws.Range("J15").Select()
